# Apply updated crypto price/volume data to Sheet1, matching the source commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain (non-ambiguous) string -- safe to set directly.
$plainUpdates = @{
    "D2" = "27.218.81"
    "E2" = "  +0.13%  "
    "D3" = "1.905.23"
    "E3" = "  +0.06%  "
    "E4" = "  +0.27%  "
    "E5" = "  +0.47%  "
    "E6" = "  +0.27%  "
    "E7" = "  +0.32%  "
    "E8" = "  +1.19%  "
    "E9" = "  +0.75%  "
    "E10" = "  +2.29%  "
    "E11" = "  +0.53%  "
    "E12" = "  -3.77%  "
    "E13" = "  +0.98%  "
    "E14" = "  +1.49%  "
    "D15" = "1.791.97"
    "E15" = "  -5.87%  "
    "E16" = "  +0.21%  "
    "E17" = "  +0.76%  "
    "E18" = "  +1.11%  "
    "D20" = "27.252.83"
    "E20" = "  +0.12%  "
    "E21" = "  +1.19%  "
    "E22" = "  +1.97%  "
    "E23" = "  +0.81%  "
    "E24" = "  +3.37%  "
    "E25" = "  +1.70%  "
    "E26" = "  +0.55%  "
    "E27" = "  -0.50%  "
    "E28" = "  +1.88%  "
    "E29" = "  +0.68%  "
    "E31" = "  -0.31%  "
    "E32" = "  -0.39%  "
    "E33" = "  +0.11%  "
    "E34" = "  -0.93%  "
    "E35" = "  +0.42%  "
    "E36" = "  +0.19%  "
    "E37" = "  +3.51%  "
    "E38" = "  -0.04%  "
    "E39" = "  +0.35%  "
    "E40" = "  +1.08%  "
    "E41" = "  +0.19%  "
    "E42" = "  -0.60%  "
    "E43" = "  -0.44%  "
    "E44" = "  +0.25%  "
    "E45" = "  +1.23%  "
    "B46" = "EnergySwap"
    "C46" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E46" = "  -0.19%  "
    "B47" = "PaxDollar"
    "C47" = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "E47" = "  +0.29%  "
    "E48" = "  +1.81%  "
    "E49" = "  +3.00%  "
    "E50" = "  +0.58%  "
    "E51" = "  +0.41%  "
}

# Cells whose new text looks like a number (e.g. "307.76") -- Excel would silently
# coerce these to numeric values unless the cell is forced to Text format first.
# The original number format/style is restored afterwards so no stray styling is left behind.
$textForcedUpdates = @{
    "D5" = "307.76"
    "D7" = "0.5254"
    "D8" = "0.3811"
    "D9" = "0.07311"
    "D11" = "0.9056"
    "D12" = "0.08090"
    "D13" = "95.82"
    "D14" = "5.366"
    "D18" = "14.75"
    "D21" = "5.128"
    "D22" = "10.83"
    "D23" = "6.486"
    "D24" = "2.364"
    "D25" = "149.82"
    "D27" = "1.743"
    "D28" = "117.10"
    "D29" = "4.848"
    "D30" = "4.888"
    "D31" = "0.09241"
    "D32" = "0.8063"
    "D33" = "0.05072"
    "D34" = "1.231"
    "D35" = "2.979"
    "D36" = "3.392"
    "D37" = "2.701"
    "D38" = "0.5730"
    "D39" = "0.01995"
    "D40" = "1.087"
    "D41" = "9.010"
    "D42" = "6.610"
    "D43" = "116.60"
    "D45" = "0.4912"
    "D46" = "10.19"
    "D47" = "1.003"
    "D50" = "64.33"
    "D51" = "0.05963"
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

foreach ($addr in $textForcedUpdates.Keys) {
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$addr]
    $cell.Style = $originalStyle
}
